$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.479.25"
$ws.Range("E2").Value = "  +0.44%  "
$ws.Range("D3").Value = "1.571.59"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "290.95"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("E7").Value = "  -1.42%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.98"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.86%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3377"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.56%  "
$ws.Range("E10").Value = "  +1.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07533"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.46%  "
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("E13").Value = "  +0.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.018"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.94%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.962"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.64%  "
$ws.Range("D16").Value = "1.572.44"
$ws.Range("E16").Value = "  +0.38%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001120"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.62%  "
$ws.Range("E18").Value = "  +0.77%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06769"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.28%  "
$ws.Range("E20").Value = "  -0.08%  "
$ws.Range("E21").Value = "  +2.29%  "
$ws.Range("E22").Value = "  -0.89%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.23"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.52%  "
$ws.Range("D24").Value = "22.483.41"
$ws.Range("E24").Value = "  +0.51%  "
$ws.Range("E25").Value = "  -0.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.615"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.49%  "
$ws.Range("E27").Value = "  -0.73%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "149.04"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.078"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.86%  "
$ws.Range("E30").Value = "  -0.39%  "
$ws.Range("D31").Value = "1.749.83"
$ws.Range("E31").Value = "  +0.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.071"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.197"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.014"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.28%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.800"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.65%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08347"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.28%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02476"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.34%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.360"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.13%  "
$ws.Range("E39").Value = "  +0.53%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06545"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.42%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.429"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.27%  "
$ws.Range("E42").Value = "  +0.36%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6216"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.37%  "
$ws.Range("E44").Value = "  +0.40%  "
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.806"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5845"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.18"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.63%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.070"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.229"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.67%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07310"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.46%  "
